$wb = $excel.ActiveWorkbook

# --- Update the "testCitizen" sheet values (A1:A8, B1:B8) ---
$ws2 = $wb.Worksheets.Item("testCitizen")

$colA = @("ulais1145","ulais1146","ulais1147","ulais1148","ulais1149","ulais1150","ulais1151","ulais1152")
$colB = @("urbs13","urbs14","urbs15","urbs16","urbs17","urbs18","urbs19","urbs20")

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 2).Value = $colB[$i]
}

# --- Update selection on testCitizen sheet and make it the active sheet/tab ---
$ws2.Range("D8").Select()
$ws2.Activate()

$wb.Save()
